# Add a new worksheet "Wheel" at the end of the workbook and populate it
# with vehicle data (Maker / Year / Model / Trim), matching the edit made
# by the workbook author ("new code by shambhu").

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the current last sheet so it lands at the end
# of the tab strip (Worksheets.Add() defaults to inserting before the
# active sheet, so we explicitly anchor it After the last existing sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Wheel"

# Header row
$ws.Range("A1").Value = "Maker"
$ws.Range("B1").Value = "Year"
$ws.Range("C1").Value = "Model"
$ws.Range("D1").Value = "Trim"

# Data rows (entered column-by-column, right to left, as the original
# author did, to match the recorded shared-string insertion order)
$ws.Range("D2").Value = "Grip"
$ws.Range("C2").Value = "AZ05"
$ws.Range("B2").Value = 2020
$ws.Range("A2").Value = "Apollo"

$ws.Range("D3").Value = "Grip"
$ws.Range("C3").Value = "AZ05"
$ws.Range("B3").Value = 2020
$ws.Range("A3").Value = "Apollo"

# Leave the cursor where the author left it before saving
$ws.Range("I6").Select()
